$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for rows 2-51 to the
# refreshed values from this run of the crypto-price scraper.
#
# D-column values that look like plain decimal numbers (e.g. "130.11") are
# written with a leading apostrophe. Excel treats a leading apostrophe as a
# "force text" prefix, so the cell keeps the exact original string (and the
# %-formatted / multi-dot "price" strings such as "61.974.61" are otherwise
# already left alone as text) instead of being auto-converted to a float and
# losing precision / its text formatting.

$ws.Range("D2").Value = "61.974.61"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'410.37"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").Value = "'130.11"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.635"
$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -2.40%  "

$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("D11").Value = "'43.69"

$ws.Range("D12").Value = "'0.0000223"
$ws.Range("E12").Value = "  +14.18%  "

$ws.Range("D13").Value = "'9.36"
$ws.Range("E13").Value = "  +4.88%  "

$ws.Range("D14").Value = "3.973.19"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "'21.22"
$ws.Range("E16").Value = "  +3.07%  "

$ws.Range("D17").Value = "3.421.99"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "'12.39"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").Value = "'1.09"
$ws.Range("E19").Value = "  +2.75%  "

$ws.Range("D20").Value = "61.970.70"
$ws.Range("E20").Value = "  -0.43%  "

$ws.Range("D21").Value = "'510.28"
$ws.Range("E21").Value = "  +25.62%  "

$ws.Range("D22").Value = "'92.93"
$ws.Range("E22").Value = "  +2.25%  "

$ws.Range("D23").Value = "'3.31"
$ws.Range("E23").Value = "  +3.63%  "

$ws.Range("E24").Value = "  +0.34%  "

$ws.Range("D25").Value = "'3.34"
$ws.Range("E25").Value = "  +3.01%  "

$ws.Range("D26").Value = "'35.10"
$ws.Range("E26").Value = "  +6.52%  "

$ws.Range("D27").Value = "'9.28"
$ws.Range("E27").Value = "  +9.09%  "

$ws.Range("D28").Value = "'7.68"
$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("E29").Value = "  +2.55%  "

$ws.Range("D30").Value = "'2.69"
$ws.Range("E30").Value = "  -2.03%  "

$ws.Range("E31").Value = "  -1.65%  "

$ws.Range("E32").Value = "  -2.64%  "

$ws.Range("D33").Value = "'41.98"
$ws.Range("E33").Value = "  -4.92%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "'58.88"
$ws.Range("E35").Value = "  +12.28%  "

$ws.Range("E36").Value = "  +1.20%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("E38").Value = "  +4.95%  "

$ws.Range("D39").Value = "'3.46"
$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("E40").Value = "  +18.44%  "

$ws.Range("E41").Value = "  +5.52%  "

$ws.Range("E42").Value = "  +1.15%  "

$ws.Range("D43").Value = "'2.12"
$ws.Range("E43").Value = "  +6.74%  "

$ws.Range("D44").Value = "'0.319"
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("D45").Value = "'4.33"
$ws.Range("E45").Value = "  +7.18%  "

$ws.Range("D46").Value = "'2.38"
$ws.Range("E46").Value = "  +23.00%  "

$ws.Range("D47").Value = "'16.64"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D48").Value = "'121.19"
$ws.Range("E48").Value = "  +23.95%  "

$ws.Range("D49").Value = "'22.96"
$ws.Range("E49").Value = "  +3.00%  "

$ws.Range("E50").Value = "  +18.50%  "

$ws.Range("D51").Value = "2.145.46"
$ws.Range("E51").Value = "  +0.96%  "
